# Modelo_de_Estimacion.xlsx update
# "Se actualizo el modelo de estimacion"
#
# The "MODELO ESTIMACION" sheet (first sheet) has its estimated-hours
# column (H) revised: most task rows double their hour estimate, while the
# very first task (row 2) quadruples it. Columns J, K, L, M, N are formulas
# that recompute automatically. The two summary cells at the bottom of
# column N (rows 49/50) swap which formula they hold.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "Horas estimadas" (column H) values, keyed by row number.
$hoursByRow = [ordered]@{
    2  = 16
    3  = 8
    6  = 6
    7  = 10
    8  = 16
    9  = 20
    10 = 6
    11 = 6
    12 = 6
    13 = 6
    14 = 2
    15 = 12
    16 = 8
    17 = 10
    18 = 2
    19 = 6
    20 = 8
    21 = 8
    22 = 2
    23 = 8
    24 = 2
    25 = 10
    26 = 10
    27 = 10
    28 = 10
    29 = 30
    30 = 4
    31 = 12
    32 = 8
    33 = 8
    34 = 2
    35 = 2
    36 = 8
    37 = 8
    38 = 2
    39 = 8
    40 = 8
    41 = 8
    42 = 8
    43 = 2
    44 = 8
    45 = 2
    46 = 8
    47 = 8
}

foreach ($row in $hoursByRow.Keys) {
    $ws.Cells.Item($row, 8).Value = $hoursByRow[$row]
}

# Bottom summary block: N49 and N50 swap formulas.
$ws.Range("N49").Formula = "=N48*1.68"
$ws.Range("N50").Formula = "=L48-N48"

# Restore the selection / active cell recorded in the sheet view.
$ws.Activate() | Out-Null
$ws.Range("I45").Select() | Out-Null
